# Apply scheduled market-price refresh to the Leve profit calculations
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 2100
$ws.Range("I18").Value = 2100
$ws.Range("K18").Value = 2100
$ws.Range("M18").Value = -1816

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2503.8064
$ws.Range("I137").Value = 1258.6666
$ws.Range("K137").Value = 3775.9998
$ws.Range("M137").Value = -1225.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 606.25
$ws.Range("I5").Value = 537.5
$ws.Range("J5").Value = 675
$ws.Range("K5").Value = 537.5
$ws.Range("L5").Value = 675
$ws.Range("M5").Value = -425.5
$ws.Range("N5").Value = -899

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 6385.9688
$ws.Range("I32").Value = 3726.8215
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 3726.8215
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -3439.8215
$ws.Range("N32").Value = -25574

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4750.8335
$ws.Range("J61").Value = 4847.8
$ws.Range("L61").Value = 4847.8
$ws.Range("N61").Value = -5271.8

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 5683.846
$ws.Range("I63").Value = 2659.8
$ws.Range("J63").Value = 7573.875
$ws.Range("K63").Value = 2659.8
$ws.Range("L63").Value = 7573.875
$ws.Range("M63").Value = -1973.8
$ws.Range("N63").Value = -8945.875

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 5683.846
$ws.Range("I66").Value = 2659.8
$ws.Range("J66").Value = 7573.875
$ws.Range("K66").Value = 13299
$ws.Range("L66").Value = 37869.375
$ws.Range("M66").Value = -9867
$ws.Range("N66").Value = -44733.375

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4750.8335
$ws.Range("J136").Value = 4847.8
$ws.Range("L136").Value = 14543.4
$ws.Range("N136").Value = -19643.4

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 606.25
$ws.Range("I4").Value = 537.5
$ws.Range("J4").Value = 675
$ws.Range("K4").Value = 537.5
$ws.Range("L4").Value = 675
$ws.Range("M4").Value = -422.5
$ws.Range("N4").Value = -905

$ws = $wb.Worksheets.Item("CRP")
# Row 13: Compulsory Conjury / Maple Cane
$ws.Range("H13").Value = 24998.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 24998.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 24998.5
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -25276.5

# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 471
$ws.Range("I22").Value = 234.5
$ws.Range("J22").Value = 660.2
$ws.Range("K22").Value = 234.5
$ws.Range("L22").Value = 660.2
$ws.Range("M22").Value = 115.5
$ws.Range("N22").Value = -1360.2

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3602.5312
$ws.Range("I31").Value = 1765.65
$ws.Range("K31").Value = 1765.65
$ws.Range("M31").Value = -1470.65

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3602.5312
$ws.Range("I34").Value = 1765.65
$ws.Range("K34").Value = 1765.65
$ws.Range("M34").Value = -1563.65

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 4488.1665
$ws.Range("I58").Value = 3852.6365
$ws.Range("K58").Value = 3852.6365
$ws.Range("M58").Value = -3649.6365

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2826.7144
$ws.Range("I132").Value = 2156.4
$ws.Range("J132").Value = 4502.5
$ws.Range("K132").Value = 6469.200000000001
$ws.Range("L132").Value = 13507.5
$ws.Range("M132").Value = -3939.200000000001
$ws.Range("N132").Value = -18567.5

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 4488.1665
$ws.Range("I136").Value = 3852.6365
$ws.Range("K136").Value = 11557.9095
$ws.Range("M136").Value = -9007.9095

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 479.82608
$ws.Range("J5").Value = 456.93332
$ws.Range("L5").Value = 1370.79996
$ws.Range("N5").Value = -1594.79996

# Row 56: Culture Club / Crowned Pie
$ws.Range("H56").Value = 12399
$ws.Range("I56").Value = 12399
$ws.Range("K56").Value = 12399
$ws.Range("M56").Value = -11869

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 625.0714
$ws.Range("I107").Value = 755
$ws.Range("J107").Value = 615.0769
$ws.Range("K107").Value = 2265
$ws.Range("L107").Value = 1845.2307
$ws.Range("M107").Value = -345
$ws.Range("N107").Value = -5685.2307

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 6343
$ws.Range("I134").Value = 4514.5
$ws.Range("K134").Value = 13543.5
$ws.Range("M134").Value = -8473.5

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 479.82608
$ws.Range("J135").Value = 456.93332
$ws.Range("L135").Value = 4112.39988
$ws.Range("N135").Value = -9182.399880000001

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 1042.591
$ws.Range("I107").Value = 1212
$ws.Range("J107").Value = 925.3077
$ws.Range("K107").Value = 1212
$ws.Range("L107").Value = 925.3077
$ws.Range("M107").Value = 708
$ws.Range("N107").Value = -4765.3077

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2959.6667
$ws.Range("I7").Value = 2959.6667
$ws.Range("K7").Value = 2959.6667
$ws.Range("M7").Value = -2847.6667

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2478.3333
$ws.Range("I40").Value = 1790
$ws.Range("K40").Value = 1790
$ws.Range("M40").Value = -1654

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 3393.8696
$ws.Range("I82").Value = 3509
$ws.Range("J82").Value = 2979.4
$ws.Range("K82").Value = 3509
$ws.Range("L82").Value = 2979.4
$ws.Range("M82").Value = -3148
$ws.Range("N82").Value = -3701.4

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 3393.8696
$ws.Range("I85").Value = 3509
$ws.Range("J85").Value = 2979.4
$ws.Range("K85").Value = 3509
$ws.Range("L85").Value = 2979.4
$ws.Range("M85").Value = -2261
$ws.Range("N85").Value = -5475.4

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2959.6667
$ws.Range("I126").Value = 2959.6667
$ws.Range("K126").Value = 8879.000100000001
$ws.Range("M126").Value = -6409.000100000001

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 7653.3076
$ws.Range("I62").Value = 5500
$ws.Range("J62").Value = 8044.8184
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 8044.8184
$ws.Range("M62").Value = -4876
$ws.Range("N62").Value = -9292.8184

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 7653.3076
$ws.Range("I65").Value = 5500
$ws.Range("J65").Value = 8044.8184
$ws.Range("K65").Value = 27500
$ws.Range("L65").Value = 40224.092
$ws.Range("M65").Value = -24380
$ws.Range("N65").Value = -46464.092

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1585.25
$ws.Range("I132").Value = 1454.5714
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4363.7142
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1833.7142
$ws.Range("N132").Value = -12560
